$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.99
$wsSummary.Range("B4").Value = -0.02
$wsSummary.Range("B5").Value = 0
$wsSummary.Range("B6").Value = 103
$wsSummary.Range("B7").Value = 39
$wsSummary.Range("B9").Value = 37.86

# ---- Sheet: Strategy Status ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.98999999999999
$wsStatus.Range("D4").Value = 103
$wsStatus.Range("E4").Value = -0.02
$wsStatus.Range("F4").Value = -0.01
$wsStatus.Range("G4").Value = 37.86

# ---- New trade row data (Trade #103) ----
$tradeNum = 103
$date = "2026-02-17"
$time = "15:58:15"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.92
$exitPrice = 0.95
$status = "CLOSED"
$pnlPct = 3.2609
$pnlDollar = 0.03
$capitalAfter = 99.98999999999999
$entrySlip = 0
$exitSlip = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.14

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 104

    $ws.Cells.Item($row, 1).Value = $tradeNum
    # Force column B to text so the ISO date string isn't auto-converted
    # into a date serial number by Excel's smart entry parsing.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $date
    $ws.Cells.Item($row, 3).Value = $time
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = $exitPrice
    $ws.Cells.Item($row, 8).Value = $status
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlip
    $ws.Cells.Item($row, 13).Value = $exitSlip
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}
